$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country re-ordering (rows whose ranking changed position) ---
# Jordania/Birmania/Serbia block: Birmania overtook Serbia, so the two
# country labels on rows 79/80 swap (data below is rewritten afterwards).
$ws.Range("A79").Value = "Birmania"
$ws.Range("A80").Value = "Serbia"

# Aruba/Islandia/Estonia/Mayotte block: Mayotte overtook Islandia and
# Estonia, so labels on rows 141/142/143 shift down one slot.
$ws.Range("A141").Value = "Mayotte"
$ws.Range("A142").Value = "Islandia"
$ws.Range("A143").Value = "Estonia"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 8390547
$ws.Range("C4").Value = 2748
$ws.Range("D4").Value = 5458024
$ws.Range("E4").Value = 2707762
$ws.Range("G4").Value = 31
$ws.Range("H4").Value = 224761

$ws.Range("B5").Value = 7555776
$ws.Range("C5").Value = 7538
$ws.Range("D5").Value = 6667565
$ws.Range("E5").Value = 773497
$ws.Range("G5").Value = 72
$ws.Range("H5").Value = 114714

$ws.Range("B21").Value = 370248
$ws.Range("C21").Value = 3267
$ws.Range("E21").Value = 68473
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 9875

$ws.Range("B49").Value = 101599
$ws.Range("C49").Value = 239
$ws.Range("D49").Value = 91032
$ws.Range("E49").Value = 7026
$ws.Range("G49").Value = 11
$ws.Range("H49").Value = 3541

$ws.Range("B51").Value = 93127
$ws.Range("C51").Value = 471
$ws.Range("D51").Value = 85941
$ws.Range("E51").Value = 5512
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 1674

$ws.Range("E57").Value = 27621
$ws.Range("G57").Value = 15
$ws.Range("H57").Value = 2138

$ws.Range("E58").Value = 3281
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 301

$ws.Range("B59").Value = 67302
$ws.Range("C59").Value = 252
$ws.Range("D59").Value = 48493
$ws.Range("E59").Value = 17209
$ws.Range("G59").Value = 16
$ws.Range("H59").Value = 1600

$ws.Range("B79").Value = 37205
$ws.Range("C79").Value = 1180
$ws.Range("D79").Value = 17568
$ws.Range("E79").Value = 18723
$ws.Range("G79").Value = 34
$ws.Range("H79").Value = 914

$ws.Range("B80").Value = 36282
$ws.Range("C80").Value = 122
$ws.Range("D80").Value = 31536
$ws.Range("E80").Value = 3968
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 778

$ws.Range("B97").Value = 16539
$ws.Range("C97").Value = 83
$ws.Range("E97").Value = 4398

$ws.Range("B98").Value = 15897
$ws.Range("C98").Value = 44
$ws.Range("D98").Value = 15031
$ws.Range("E98").Value = 520

$ws.Range("B111").Value = 10533
$ws.Range("C111").Value = 40
$ws.Range("D111").Value = 9563
$ws.Range("E111").Value = 890

$ws.Range("B115").Value = 8321
$ws.Range("C115").Value = 47
$ws.Range("D115").Value = 3951
$ws.Range("E115").Value = 4197
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 173

$ws.Range("B122").Value = 6258
$ws.Range("C122").Value = 38
$ws.Range("D122").Value = 5780
$ws.Range("E122").Value = 351
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 127

$ws.Range("B126").Value = 5625
$ws.Range("C126").Value = 87
$ws.Range("E126").Value = 2172

$ws.Range("B141").Value = 4159
$ws.Range("C141").Value = 129
$ws.Range("D141").Value = 2964
$ws.Range("E141").Value = 1152
$ws.Range("H141").Value = 43

$ws.Range("B142").Value = 4101
$ws.Range("C142").Value = 46
$ws.Range("D142").Value = 2856
$ws.Range("E142").Value = 1234
$ws.Range("H142").Value = 11

$ws.Range("B143").Value = 4085
$ws.Range("C143").Value = 7
$ws.Range("D143").Value = 3229
$ws.Range("E143").Value = 788
$ws.Range("H143").Value = 68

$ws.Range("B177").Value = 549
$ws.Range("C177").Value = 7
$ws.Range("D177").Value = 51

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 16:58"
